$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALERTS")

$row = 5
$cells = $ws.Range("A" + $row + ":F" + $row)
$cells.NumberFormat = "@"

$ws.Range("A5").Value = "2026-01-31"
$ws.Range("B5").Value = "21:49:48"
$ws.Range("C5").Value = "21:00"
$ws.Range("D5").Value = "Living Room"
$ws.Range("E5").Value = "CRITICAL"
$ws.Range("F5").Value = "FALL_DETECTED"
